$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report period) ---
$ws.Range("A8").Value = "Volume 33   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/5/2026  Through  1/11/2026"

# --- Simple value updates (style/number format unchanged) ---
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("L15").Value = 0
$ws.Range("C16").Value = 11
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = 22.222222222222
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 34
$ws.Range("H16").Value = -17.647058823529
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 15
$ws.Range("K16").Value = -13.333333333333
$ws.Range("L16").Value = 8.333333333333
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = -72.916666666666
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 61
$ws.Range("G17").Value = 51
$ws.Range("H17").Value = 19.607843137254
$ws.Range("I17").Value = 23
$ws.Range("J17").Value = 15
$ws.Range("K17").Value = 53.333333333333
$ws.Range("L17").Value = 15
$ws.Range("M17").Value = 187.5
$ws.Range("N17").Value = -25.806451612903
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -16.666666666666
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -38.461538461538
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 8
$ws.Range("K18").Value = -25
$ws.Range("L18").Value = -14.285714285714
$ws.Range("N18").Value = -76
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = -11.627906976744
$ws.Range("I19").Value = 13
$ws.Range("J19").Value = 20
$ws.Range("K19").Value = -35
$ws.Range("L19").Value = -38.095238095238
$ws.Range("M19").Value = 85.714285714285
$ws.Range("N19").Value = -18.75
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 14
$ws.Range("H20").Value = 7.692307692307
$ws.Range("I20").Value = 4
$ws.Range("L20").Value = -50
$ws.Range("M20").Value = 33.333333333333
$ws.Range("N20").Value = -75
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 160
$ws.Range("G21").Value = 171
$ws.Range("H21").Value = -6.432748538011
$ws.Range("I21").Value = 61
$ws.Range("J21").Value = 62
$ws.Range("K21").Value = -1.612903225806
$ws.Range("L21").Value = -11.594202898550
$ws.Range("M21").Value = 103.333333333333
$ws.Range("N21").Value = -55.147058823529
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -80
$ws.Range("F23").Value = 15
$ws.Range("H23").Value = -21.052631578947
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 7
$ws.Range("K23").Value = -42.857142857142
$ws.Range("L23").Value = -76.470588235294
$ws.Range("M23").Value = -71.428571428571
$ws.Range("F24").Value = 69
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = -22.471910112359
$ws.Range("I24").Value = 17
$ws.Range("J24").Value = 33
$ws.Range("K24").Value = -48.484848484848
$ws.Range("L24").Value = -29.166666666666
$ws.Range("M24").Value = -10.526315789473
$ws.Range("F25").Value = 11
$ws.Range("H25").Value = 22.222222222222
$ws.Range("I25").Value = 3
$ws.Range("J25").Value = 4
$ws.Range("K25").Value = -25
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 18
$ws.Range("E26").Value = -27.777777777777
$ws.Range("F26").Value = 69
$ws.Range("G26").Value = 65
$ws.Range("H26").Value = 6.153846153846
$ws.Range("I26").Value = 26
$ws.Range("J26").Value = 24
$ws.Range("K26").Value = 8.333333333333
$ws.Range("L26").Value = 62.5
$ws.Range("M26").Value = 23.809523809523
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("L27").Value = -50
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 350
$ws.Range("I28").Value = 2
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 2
$ws.Range("N29").Value = 100
$ws.Range("F30").Value = 2
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 2
$ws.Range("N30").Value = 100
$ws.Range("J44").Value = 583
$ws.Range("K44").Value = 176.303317535545
$ws.Range("L44").Value = 177.619047619048
$ws.Range("M44").Value = 60.606060606060
$ws.Range("N44").Value = 52.617801047120
$ws.Range("J46").Value = 2319
$ws.Range("K46").Value = 40.972644376899
$ws.Range("L46").Value = 13.899803536345
$ws.Range("M46").Value = -46.132404181184
$ws.Range("N46").Value = -47.581374321880

# --- Cells moving from text placeholder ("0"/"***.*") to real numbers ---
# NumberFormat must be (re)applied so the cell style matches the numeric style
# used elsewhere in the table (#,##0 for counts, #,##0.0;"-"#,##0.0 for % chg).
$ws.Range("C14").Value = 1
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("F14").Value = 1
$ws.Range("F14").NumberFormat = "#,##0"
$ws.Range("I14").Value = 1
$ws.Range("I14").NumberFormat = "#,##0"
$ws.Range("M14").Value = 0
$ws.Range("M14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = 0
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I15").Value = 1
$ws.Range("I15").NumberFormat = "#,##0"
$ws.Range("J15").Value = 1
$ws.Range("J15").NumberFormat = "#,##0"
$ws.Range("K15").Value = 0
$ws.Range("K15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M18").Value = 500
$ws.Range("M18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J20").Value = 3
$ws.Range("J20").NumberFormat = "#,##0"
$ws.Range("K20").Value = 33.333333333333
$ws.Range("K20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I27").Value = 1
$ws.Range("I27").NumberFormat = "#,##0"
$ws.Range("J27").Value = 1
$ws.Range("J27").NumberFormat = "#,##0"
$ws.Range("K27").Value = 0
$ws.Range("K27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = 0
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J28").Value = 1
$ws.Range("J28").NumberFormat = "#,##0"
$ws.Range("K28").Value = 100
$ws.Range("K28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L28").Value = 100
$ws.Range("L28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L29").Value = 100
$ws.Range("L29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L30").Value = 100
$ws.Range("L30").NumberFormat = "#,##0.0;""-""#,##0.0"

# --- Cell moving from a number back to a text placeholder ---
# Copy directly from a sibling cell that already holds the exact target
# text/format ("0", right-aligned, General) so the style is reused exactly.
$ws.Range("D22").Copy($ws.Range("C22"))
